$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "last_edited_time" column (D) timestamps recorded for this Notion export
# need to move forward by 11 minutes overall, with the group boundary for
# rows 18-21 shifting from the first group into the second, and for rows
# 120-127 shifting from the third group into the fourth group.

# Rows 2-17 -> 2024-08-03T03:28:00.000Z
$ws.Range("D2:D17").Value = "2024-08-03T03:28:00.000Z"

# Rows 18-76 -> 2024-08-03T03:29:00.000Z
$ws.Range("D18:D76").Value = "2024-08-03T03:29:00.000Z"

# Rows 77-119 -> 2024-08-03T03:30:00.000Z
$ws.Range("D77:D119").Value = "2024-08-03T03:30:00.000Z"

# Rows 120-140 -> 2024-08-03T03:31:00.000Z
$ws.Range("D120:D140").Value = "2024-08-03T03:31:00.000Z"
